$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 4
    4  = 6
    5  = 12
    6  = 6
    7  = 5
    8  = 3
    9  = 2
    10 = 6
    11 = 5
    12 = 5
    13 = 6
    14 = 4
    15 = 9
    16 = 7
    17 = 5
    18 = 6
    19 = 4
    20 = 7
    21 = 7
    22 = 1
    23 = 11
    24 = 4
    25 = 10
    26 = 4
    27 = 1
    28 = 5
    29 = 7
    30 = 6
    31 = 8
    32 = 4
    33 = 1
    34 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
